$d = $word.ActiveDocument

function Set-ParaText($searchText, $newText) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like "*$searchText*") {
            $p.Range.Text = $newText
            return $true
        }
    }
    return $false
}

# 1. Remove the stray trailing line break embedded in the w:t run for the
#    "annotatable resource" paragraph.
$r1 = Set-ParaText "annotatable resource" "This is an annotatable resource in the casebook."
if (-not $r1) { throw "Could not find paragraph 1 (annotatable resource)" }
Write-Host "r1=$r1"

# 2. Collapse the embedded line breaks in the "highlighted: ..." paragraph
#    into single spaces.
$newText2 = "highlighted: content to highlight; elided: content to elide; replaced: content to replace; linked: content to link; noted: content to note; highlighted2: second highlight content;"
$r2 = Set-ParaText "highlighted: content to highlight" $newText2
if (-not $r2) { throw "Could not find paragraph 2 (highlighted: content to highlight)" }
Write-Host "r2=$r2"

# 3. Remove the stray trailing line break embedded in the w:t run for the
#    "second chapter" paragraph.
$r3 = Set-ParaText "second chapter" "This is the second chapter of the casebook."
if (-not $r3) { throw "Could not find paragraph 3 (second chapter)" }
Write-Host "r3=$r3"
